$d = $word.ActiveDocument

# --- Update the date heading ---
$d.Content.Find.Execute("2025-10-20 Monday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-10-21 Tuesday", 2)

# --- Update the practice-problem table ---
# The table has 20 rows; every 4th row (1,5,9,13,17) holds 5 data cells,
# the other rows are blank spacer rows. Cells are addressed positionally
# (Table.Cell(row, col)) so duplicate/overlapping values across cells are
# never ambiguous.
$t = $d.Tables.Item(1)

$rowUpdates = @(
    @{ Row = 1;  Values = @("18÷4=4, 2",  "91÷7=13, 0", "38÷2=19, 0", "28÷3=9, 1",  "85÷7=12, 1") },
    @{ Row = 5;  Values = @("24÷5=4, 4",  "58÷9=6, 4",  "49÷5=9, 4",  "46÷7=6, 4",  "25÷6=4, 1") },
    @{ Row = 9;  Values = @("15÷5=3, 0",  "60÷7=8, 4",  "37÷2=18, 1", "98÷8=12, 2", "49÷2=24, 1") },
    @{ Row = 13; Values = @("93÷6=15, 3", "42÷8=5, 2",  "62÷4=15, 2", "70÷5=14, 0", "48÷6=8, 0") },
    @{ Row = 17; Values = @("46÷6=7, 4",  "75÷5=15, 0", "85÷2=42, 1", "95÷2=47, 1", "42÷8=5, 2") }
)

foreach ($update in $rowUpdates) {
    $rowIndex = $update.Row
    $values = $update.Values
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}

Write-Output "done"
